$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the threshold values in column C
$ws.Range("C2").Value = 9.5
$ws.Range("C3").Value = 7.5
$ws.Range("C5").Value = 20

# Move the active selection to C3
$ws.Range("C3").Select()
